$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSetInteractionPages")

# Row 16
$ws.Range("A16").Value = "DefaultFunctionality_DragLastItemToListTop_LastItemIsAtListTop"
$ws.Range("D16").Value = 3
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").Value = 6
$ws.Range("E16").NumberFormat = "@"
$ws.Range("H16").Value = "Sortable"

# Row 17
$ws.Range("A17").Value = "ConnectLists_DragOneSortableItemFromFirstListToSecondList_SortableItemItemMovedToSecondList"
$ws.Range("D17").Value = 5
$ws.Range("D17").NumberFormat = "@"
$ws.Range("H17").Value = "Sortable"

# Update selection to match post-edit state
$ws.Range("H20").Select()

# Restore the originally active sheet so the workbook-level active tab is unchanged
$wb.Worksheets.Item("DataSetRegistrationUser").Activate()
